# Apply the target edit to the presentation:
#  1. Mark slides 8, 9, 10, 12, 13, 15, 16, 17, 18, 19, 20, 21, 22 as hidden
#     (adds show="0" to the <p:sld> element for each).
#  2. Delete slide 24 ("Comments") entirely, which also removes its
#     <p:sldId> entry from the presentation's slide list.

$p = $ppt.ActivePresentation

$hiddenSlideNumbers = @(8, 9, 10, 12, 13, 15, 16, 17, 18, 19, 20, 21, 22)

foreach ($num in $hiddenSlideNumbers) {
    $slide = $p.Slides.Item($num)
    $slide.SlideShowTransition.Hidden = $true
}

# Remove the trailing "Comments" slide (slide 24) from the deck.
$p.Slides.Item($p.Slides.Count).Delete()
